# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values for the second data
# row (row 3) on both the "zh-cn" and "de-de" worksheets, reflecting the
# newly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-02-17 04:24:29"
$wsZhCn.Range("G3").Value = "2016-02-17 04:25:15"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-02-17 04:24:39"
$wsDeDe.Range("G3").Value = "2016-02-17 04:25:33"
